# Backlog.xlsx edit
# - Column C ("Semana") on both sheets used to hold the text "Semana 01" for
#   every data row. It is changed to the plain numeric value 1.
# - On the ITI sheet the "Semana" cells also lose their shaded/fill style,
#   ending up with the same plain centered style used on the SPN sheet.
# - The active sheet/selection bookmarks change: SPN no longer is the
#   selected tab (its selection becomes C2), while ITI becomes the
#   selected tab with D24 selected.

$wb  = $excel.ActiveWorkbook
$spn = $wb.Worksheets.Item("SPN")
$iti = $wb.Worksheets.Item("ITI")

# --- SPN sheet: "Semana" column (C2:C26) text -> number 1 ---
$spn.Range("C2:C26").Value = 1

# --- ITI sheet: "Semana" column (C2:C10) text -> number 1 ---
$iti.Range("C2:C10").Value = 1
# Match the plain centered style (no fill) used elsewhere, instead of the
# shaded style the text cells used to carry.
$iti.Range("C2:C10").HorizontalAlignment = -4108

# --- View / selection state ---
# SPN keeps its grid position but is no longer the active tab; its stored
# selection moves to C2.
$spn.Range("C2").Select()

# ITI becomes the active tab with D24 selected.
$iti.Activate()
$iti.Range("D24").Select()
